# Insert a new weekly data row into the Pomelo price sheet.
# A new record (row 62) is inserted, pushing all existing rows
# from 62 downward down by one (old row 161 becomes row 162).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 62, shifting rows 62..161 down to 63..162
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new data point.
$ws.Cells.Item(62, 1).Value = 4
$ws.Cells.Item(62, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(62, 3).Value = "Los Lagos"

$ws.Cells.Item(62, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(62, 4).Value = 44495

$ws.Cells.Item(62, 5).Value = 10
$ws.Cells.Item(62, 6).Value = "Fruta"
$ws.Cells.Item(62, 7).Value = 100102
$ws.Cells.Item(62, 8).Value = "Cítricos"
$ws.Cells.Item(62, 9).Value = 100102006
$ws.Cells.Item(62, 10).Value = "Pomelo"
$ws.Cells.Item(62, 11).Value = "Start Ruby"
$ws.Cells.Item(62, 12).Value = "Primera"
$ws.Cells.Item(62, 13).Value = 240
$ws.Cells.Item(62, 14).Value = 11000
$ws.Cells.Item(62, 15).Value = 12000
$ws.Cells.Item(62, 16).Value = 11500
$ws.Cells.Item(62, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(62, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(62, 19).Value = 821
$ws.Cells.Item(62, 20).Value = 14
